$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force a number-looking string (e.g. "3.7") to be stored as text,
    # matching how the "price" column is authored elsewhere in the sheet,
    # then drop the temporary number-format style so the cell keeps the
    # workbook's default style (no explicit s="...").
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# --- Row 2: replace the "FRIES" entry with a "fries" entry (dedupe/fix) ---
$ws.Cells.Item(2,1).Value = "9fdbdc80-b998-444e-a3af-0017900c9fc2"
$ws.Cells.Item(2,2).Value = "fries"
Set-TextValue $ws.Cells.Item(2,3) "3.7"
$ws.Cells.Item(2,4).Value = "NTU"
$ws.Cells.Item(2,5).Value = "side"
$ws.Cells.Item(2,6).Value = "Hot piping fries"

# --- Row 3: replace the "3PC set meal" entry with another "fries" entry ---
$ws.Cells.Item(3,1).Value = "931fe060-d1b5-487f-9bc4-ab60bd3fa135"
$ws.Cells.Item(3,2).Value = "fries"
Set-TextValue $ws.Cells.Item(3,3) "4.3"
$ws.Cells.Item(3,4).Value = "JP"
$ws.Cells.Item(3,5).Value = "side"
$ws.Cells.Item(3,6).Value = "Hot piping fries"

# --- Row 10 (new): add a "Chicken tenders" menu item ---
$ws.Cells.Item(10,1).Value = "52805239-f4e6-494e-935e-59fd932fbb89"
$ws.Cells.Item(10,2).Value = "Chicken tenders"
Set-TextValue $ws.Cells.Item(10,3) "6.9"
$ws.Cells.Item(10,4).Value = "NTU"
$ws.Cells.Item(10,5).Value = "side"
$ws.Cells.Item(10,6).Value = "Fresh chicken"

# --- Row 11 (new blank spacer row, matches the style used by the other blank rows) ---
$ws.Cells.Item(12,1).Copy()
$ws.Cells.Item(11,1).PasteSpecial(-4122)
$ws.Cells.Item(11,1).Value = ""

# --- Re-home the three "highlighted" blank rows two rows earlier each ---
# (row 23 -> row 21, row 34 -> row 32, row 60 -> row 58), and restore the
# vacated rows back to the plain style used elsewhere in the block.
$ws.Cells.Item(23,1).Copy()
$ws.Cells.Item(21,1).PasteSpecial(-4122)

$ws.Cells.Item(22,1).Copy()
$ws.Cells.Item(23,1).PasteSpecial(-4122)

$ws.Cells.Item(34,1).Copy()
$ws.Cells.Item(32,1).PasteSpecial(-4122)

$ws.Cells.Item(33,1).Copy()
$ws.Cells.Item(34,1).PasteSpecial(-4122)

$ws.Cells.Item(60,1).Copy()
$ws.Cells.Item(58,1).PasteSpecial(-4122)

$ws.Cells.Item(59,1).Copy()
$ws.Cells.Item(60,1).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Drop the two now-unused trailing blank rows so the sheet ends at 72 ---
$ws.Rows(73).Delete()
$ws.Rows(73).Delete()
